$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new row 71 is (apart from a few fields) a near-duplicate of row 70,
# and in particular column A on row 70 currently holds the phone number
# "71717173" stored as text - which is exactly what A71 needs to be too.
# Copy row 70 as values first (before we change A70's type below) so that
# A71 inherits the original text data type instead of being
# re-interpreted as a number by a plain .Value assignment.
$ws.Range("A70:J70").Copy()
$ws.Range("A71").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# Row 70 / A70: was stored as text "71717173", now becomes the numeric
# value 71717173.
$ws.Range("A70").Value = 71717173

# Fill in row 71's fields that differ from the copied row 70 template.
$ws.Range("C71").Value = "Cash"
$ws.Range("D71").Value = "2025-08-20T08:25:30"
$ws.Range("E71").Value = 351
$ws.Range("G71").Value = 298.35
$ws.Range("H71").Value = 52.65
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
